# ---------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored
# as TEXT (inline/shared string) even when it looks like a number
# or a date (e.g. "2", "99", "4/20/2023"), then strip the leftover
# "@" number-format override so the cell ends up with the same
# (inherited / default) style it would have had if the text had
# simply been typed into an already-General-formatted cell.
# ---------------------------------------------------------------
function Set-TextValue {
    param($cell, $value, $scratch)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $scratch.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats -> only pulls in the (clean) format
}

$wb = $excel.ActiveWorkbook

# ===================================================================
# Sheet "Sheet1"  (dimension A1:F23 -> A1:F21)
# ===================================================================
$ws1 = $wb.Worksheets.Item("Sheet1")
$scratch1 = $ws1.Cells.Item(500, 26)   # untouched scratch cell (Z500)

function Set-Row1 {
    param($r, $a, $b, $c, $d, $e, $f)
    $ws1.Cells.Item($r, 1).Value = $a
    Set-TextValue $ws1.Cells.Item($r, 2) $b $scratch1
    Set-TextValue $ws1.Cells.Item($r, 3) $c $scratch1
    Set-TextValue $ws1.Cells.Item($r, 4) $d $scratch1
    $ws1.Cells.Item($r, 5).Value = $e
    $ws1.Cells.Item($r, 6).Value = $f
}

Set-Row1 16 "Joe"               "4/20/2023"  "2" "99" "FB" "Strike looking"
Set-Row1 17 "Joe"               "4/20/2023"  "4" "88" "SL" "Foul Ball"
Set-Row1 18 "Joe"               "4/20/2023"  "6" "99" "FB" "Strikeout looking"
Set-Row1 19 "Woody"             "03/30/2023" "1" "98" "FB" "Strike looking"
Set-Row1 20 "Andrew Armstrong"  "09/23/2023" "1" "98" "CB" "Strike looking"
Set-Row1 21 "Andrew Armstrong"  "09/22/2023" "1" "23" "FB" "Strike swing & miss"

$ws1.Rows("22:23").Delete()

# ===================================================================
# Sheet "pitch breakdown"  (dimension A1:I23 -> A1:I21)
# ===================================================================
$ws2 = $wb.Worksheets.Item("pitch breakdown")
$scratch2 = $ws2.Cells.Item(500, 26)   # untouched scratch cell (Z500)

function Set-Row2 {
    param($r, $a, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws2.Cells.Item($r, 1).Value = $a
    Set-TextValue $ws2.Cells.Item($r, 2) $b $scratch2
    Set-TextValue $ws2.Cells.Item($r, 3) $c $scratch2
    Set-TextValue $ws2.Cells.Item($r, 4) $d $scratch2
    $ws2.Cells.Item($r, 5).Value = $e
    $ws2.Cells.Item($r, 6).Value = $f
    $ws2.Cells.Item($r, 7).Value = $g
    $ws2.Cells.Item($r, 8).Value = $h
    $ws2.Cells.Item($r, 9).Value = $i
}

Set-Row2 16 "Joe"               "4/20/2023"  "2" "99" "FB" "Strike looking"     "Strike" "No swing"         "nothing"
Set-Row2 17 "Joe"               "4/20/2023"  "4" "88" "SL" "Foul Ball"          "Strike" "Swing contact"    "nothing"
Set-Row2 18 "Joe"               "4/20/2023"  "6" "99" "FB" "Strikeout looking"  "Strike" "No swing"         "not free base"
Set-Row2 19 "Woody"             "03/30/2023" "1" "98" "FB" "Strike looking"    "Strike" "No swing"         "nothing"
Set-Row2 20 "Andrew Armstrong"  "09/23/2023" "1" "98" "CB" "Strike looking"    "Strike" "No swing"         "nothing"
Set-Row2 21 "Andrew Armstrong"  "09/22/2023" "1" "23" "FB" "Strike swing & miss" "Strike" "Swing no contact" "nothing"

$ws2.Rows("22:23").Delete()

# ===================================================================
# Sheet "pitcher breakdown"  (dimension A1:I5 -> A1:I6)
# ===================================================================
$ws3 = $wb.Worksheets.Item("pitcher breakdown")

$ws3.Cells.Item(5, 1).Value = "Joe"
$ws3.Cells.Item(5, 2).Value = 99
$ws3.Cells.Item(5, 3).Value = 99
$ws3.Cells.Item(5, 4).Value = 1
$ws3.Cells.Item(5, 5).Value = 0
$ws3.Cells.Item(5, 6).Value = 0.7
$ws3.Cells.Item(5, 7).Value = 1
$ws3.Cells.Item(5, 8).Value = 0
$ws3.Cells.Item(5, 9).Value = 0

$ws3.Cells.Item(6, 1).Value = "Andrew Armstrong"
$ws3.Cells.Item(6, 2).Value = 23
$ws3.Cells.Item(6, 3).Value = 23
$ws3.Cells.Item(6, 4).Value = 1
$ws3.Cells.Item(6, 5).Value = 1
$ws3.Cells.Item(6, 6).Value = 1
$ws3.Cells.Item(6, 7).Value = 1
$ws3.Cells.Item(6, 8).Value = 1
$ws3.Cells.Item(6, 9).Value = 0

$excel.CutCopyMode = $false
